# first step towards #7
#
# LinkML model changes reflected in the generated Excel template:
#  - AssemblyJoinComponent + AssemblyJoin are merged into a single
#    "AssemblyFragment" sheet (sequence, left_location, right_location,
#    reverse_complemented); the separate AssemblyJoin (left/right) sheet
#    is removed entirely.
#  - RepositoryIdSource / AddGeneIdSource / BenchlingUrlSource now list
#    repository_id before repository_name (column order swapped), and
#    the "addgene,genbank,benchling" dropdown validation moves from the
#    repository_name column to the repository_id column.

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) AssemblyJoinComponent -> AssemblyFragment: add left_location /
#    right_location columns (replacing the single "location" column).
# -------------------------------------------------------------------
$wsFragment = $wb.Worksheets.Item("AssemblyJoinComponent")
$wsFragment.Name = "AssemblyFragment"
$wsFragment.Range("A1").Value = "sequence"
$wsFragment.Range("B1").Value = "left_location"
$wsFragment.Range("C1").Value = "right_location"
$wsFragment.Range("D1").Value = "reverse_complemented"

# -------------------------------------------------------------------
# 2) AssemblyJoin ("left"/"right" helper sheet) is no longer part of
#    the model - drop it. Every following sheet shifts up by one.
# -------------------------------------------------------------------
[void]$wb.Worksheets.Item("AssemblyJoin").Delete()

# -------------------------------------------------------------------
# 3) RepositoryIdSource: swap repository_name/repository_id columns
#    and move the dropdown validation from A to B.
# -------------------------------------------------------------------
$wsRepo = $wb.Worksheets.Item("RepositoryIdSource")
$wsRepo.Range("A1").Value = "repository_id"
$wsRepo.Range("B1").Value = "repository_name"
[void]$wsRepo.Range("A2:A1048576").Validation.Delete()
[void]$wsRepo.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')

# -------------------------------------------------------------------
# 4) AddGeneIdSource: swap repository_name/repository_id columns
#    (columns C/D) and move the dropdown validation from C to D.
# -------------------------------------------------------------------
$wsAddGene = $wb.Worksheets.Item("AddGeneIdSource")
$wsAddGene.Range("C1").Value = "repository_id"
$wsAddGene.Range("D1").Value = "repository_name"
[void]$wsAddGene.Range("C2:C1048576").Validation.Delete()
[void]$wsAddGene.Range("D2:D1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')

# -------------------------------------------------------------------
# 5) BenchlingUrlSource: swap repository_name/repository_id columns
#    and move the dropdown validation from A to B.
# -------------------------------------------------------------------
$wsBenchling = $wb.Worksheets.Item("BenchlingUrlSource")
$wsBenchling.Range("A1").Value = "repository_id"
$wsBenchling.Range("B1").Value = "repository_name"
[void]$wsBenchling.Range("A2:A1048576").Validation.Delete()
[void]$wsBenchling.Range("B2:B1048576").Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')
